## Data for Economic Rework.xlsx -- "Factories - South America" edit
##
## 1. Highlight (fill-style) a batch of "New Factories" sheet country rows
##    that now have data entered for them (B column style 35 -> 60).
## 2. Insert a new country row (Venezuela) at row 189, pushing
##    Vietnam/Yemen/Zambia/Zimbabwe down by one row, and give it its data.
## 3. Restore/adjust sheet view (frozen-pane scroll position + selection)
##    on "New Factories" and "Support Data" sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("New Factories")

## --- 1. Mark rows that now have factory data (fill highlight s=35 -> s=60) ---
$highlightRows = @(22,36,39,40,45,50,51,52,65,86,92,99,104,108,118,128,132,147,148,149,170,176,177,187)

$ws.Range("B17").Copy() | Out-Null
foreach ($r in $highlightRows) {
    $ws.Range("B$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

## --- 2. Insert the new "Venezuela" row above Vietnam (row 189) ---
$ws.Rows.Item(189).Insert() | Out-Null

## Row-level format (matches the other "manually updated" rows, e.g. Yemen)
$ws.Range("B190:AN190").Copy() | Out-Null
$ws.Range("B189:AN189").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

## Give B189 the highlighted fill too (same as the other populated rows)
$ws.Range("B17").Copy() | Out-Null
$ws.Range("B189").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$r = 189
$ws.Range("B$r").Value2 = "Venezuela"
$ws.Range("C$r").Value2 = 168548247978
$ws.Range("D$r").Value2 = 1.5265701288490501
$ws.Range("E$r").Formula = "=D$r/100*C$r"
$ws.Range("F$r").Formula = "=ROUND(E$r/`$C`$6,0)"
$ws.Range("G$r").Formula = "=F$r*`$C`$6*`$C`$8"
$ws.Range("H$r").Formula = "=C$r-G$r"
$ws.Range("I$r").Value2 = 43.081849066103153
$ws.Range("J$r").Formula = "=I$r/(I$r+L$r)*H$r*`$F`$7"
$ws.Range("K$r").Formula = "=ROUND(J$r/`$F`$9,0)"
$ws.Range("L$r").Value2 = 46.404848256167305
$ws.Range("M$r").Formula = "=L$r/(I$r+L$r)*H$r*(2-`$F`$7)"
$ws.Range("N$r").Formula = "=ROUND(M$r/`$F`$8,0)"
$ws.Range("O$r").Formula = "=F$r"
$ws.Range("P$r").Formula = "=N$r"
$ws.Range("Q$r").Formula = "=K$r"
$ws.Range("R$r").Formula = "=O$r+P$r+Q$r"
$ws.Range("T$r").Value2 = 141632345013
$ws.Range("U$r").Formula = "=(T$r-C$r)*`$U`$11+C$r"
$ws.Range("V$r").Value2 = 0.48784414488216299
$ws.Range("W$r").Formula = "=V$r/100*U$r"
$ws.Range("X$r").Formula = "=ROUND(W$r/`$C`$6,0)"
$ws.Range("Y$r").Formula = "=X$r*`$C`$6*`$C`$8"
$ws.Range("Z$r").Formula = "=U$r-Y$r"
$ws.Range("AA$r").Value2 = 43.081849066103153
$ws.Range("AB$r").Formula = "=AA$r/(AA$r+AD$r)*Z$r*`$F`$7"
$ws.Range("AC$r").Formula = "=ROUND(AB$r/`$F`$9,0)"
$ws.Range("AD$r").Value2 = 46.404848256167305
$ws.Range("AE$r").Formula = "=AD$r/(AA$r+AD$r)*Z$r*(2-`$F`$7)"
$ws.Range("AF$r").Formula = "=ROUND(AE$r/`$F`$8,0)"
$ws.Range("AG$r").Formula = "=X$r"
$ws.Range("AH$r").Formula = "=AF$r"
$ws.Range("AI$r").Formula = "=AC$r"
$ws.Range("AJ$r").Formula = "=AG$r+AH$r+AI$r"
$ws.Range("AK$r").Formula = "=AG$r-O$r"
$ws.Range("AL$r").Formula = "=AH$r-P$r"
$ws.Range("AM$r").Formula = "=AI$r-Q$r"
$ws.Range("AN$r").Formula = "=AJ$r-R$r"

## --- 3. Sheet view / scroll-position changes ---
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("B39").Select()

$support = $wb.Worksheets.Item("Support Data")
$support.Activate()
$support.Application.ActiveWindow.ScrollColumn = 10
$support.Range("Z212").Select()

$ws.Activate()
